# =============================================================
# New PO forecast model
# Updates the three PO-analysis sheets:
#   "Weekly Quantity" - append the latest weekly PO row
#   "Monthly Trend"   - append the latest monthly PO row
#   "PO Forecast"     - re-run forecast: update historical values
#                       and extend the forecast horizon
# =============================================================

$wb = $excel.ActiveWorkbook

# --- Sheet "Weekly Quantity": append new weekly actuals row ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Cells.Item(77, 1).Value = 45662.99999999999
$wsWeekly.Cells.Item(77, 1).NumberFormat = $wsWeekly.Cells.Item(76, 1).NumberFormat
$wsWeekly.Cells.Item(77, 2).Value = 10

# --- Sheet "Monthly Trend": append new monthly actuals row ---
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Cells.Item(24, 1).Value = 45688.99999999999
$wsMonthly.Cells.Item(24, 1).NumberFormat = $wsMonthly.Cells.Item(23, 1).NumberFormat
$wsMonthly.Cells.Item(24, 2).Value = 10

# --- Sheet "PO Forecast": refreshed forecast values for existing dates ---
$wsForecast = $wb.Worksheets.Item("PO Forecast")
$wsForecast.Cells.Item(2, 2).Value = 241
$wsForecast.Cells.Item(3, 2).Value = 242
$wsForecast.Cells.Item(4, 2).Value = 245
$wsForecast.Cells.Item(5, 2).Value = 246
$wsForecast.Cells.Item(6, 2).Value = 247
$wsForecast.Cells.Item(7, 2).Value = 248
$wsForecast.Cells.Item(8, 2).Value = 249
$wsForecast.Cells.Item(9, 2).Value = 250
$wsForecast.Cells.Item(10, 2).Value = 254
$wsForecast.Cells.Item(11, 2).Value = 255
$wsForecast.Cells.Item(12, 2).Value = 256
$wsForecast.Cells.Item(13, 2).Value = 257
$wsForecast.Cells.Item(14, 2).Value = 258
$wsForecast.Cells.Item(15, 2).Value = 259
$wsForecast.Cells.Item(16, 2).Value = 260
$wsForecast.Cells.Item(17, 2).Value = 261
$wsForecast.Cells.Item(18, 2).Value = 262
$wsForecast.Cells.Item(19, 2).Value = 263
$wsForecast.Cells.Item(20, 2).Value = 264
$wsForecast.Cells.Item(21, 2).Value = 267
$wsForecast.Cells.Item(22, 2).Value = 269
$wsForecast.Cells.Item(23, 2).Value = 270
$wsForecast.Cells.Item(24, 2).Value = 272
$wsForecast.Cells.Item(25, 2).Value = 273
$wsForecast.Cells.Item(29, 2).Value = 277
$wsForecast.Cells.Item(30, 2).Value = 278
$wsForecast.Cells.Item(31, 2).Value = 279
$wsForecast.Cells.Item(32, 2).Value = 280
$wsForecast.Cells.Item(33, 2).Value = 281
$wsForecast.Cells.Item(34, 2).Value = 282
$wsForecast.Cells.Item(35, 2).Value = 283
$wsForecast.Cells.Item(36, 2).Value = 284
$wsForecast.Cells.Item(37, 2).Value = 285
$wsForecast.Cells.Item(38, 2).Value = 286
$wsForecast.Cells.Item(39, 2).Value = 287
$wsForecast.Cells.Item(40, 2).Value = 288
$wsForecast.Cells.Item(41, 2).Value = 291
$wsForecast.Cells.Item(42, 2).Value = 293
$wsForecast.Cells.Item(43, 2).Value = 295
$wsForecast.Cells.Item(44, 2).Value = 296
$wsForecast.Cells.Item(45, 2).Value = 298
$wsForecast.Cells.Item(46, 2).Value = 300
$wsForecast.Cells.Item(47, 2).Value = 301
$wsForecast.Cells.Item(48, 2).Value = 302
$wsForecast.Cells.Item(49, 2).Value = 303
$wsForecast.Cells.Item(50, 2).Value = 304
$wsForecast.Cells.Item(51, 2).Value = 305
$wsForecast.Cells.Item(52, 2).Value = 306
$wsForecast.Cells.Item(53, 2).Value = 307
$wsForecast.Cells.Item(54, 2).Value = 309
$wsForecast.Cells.Item(55, 2).Value = 310
$wsForecast.Cells.Item(56, 2).Value = 311
$wsForecast.Cells.Item(57, 2).Value = 313
$wsForecast.Cells.Item(58, 2).Value = 314
$wsForecast.Cells.Item(59, 2).Value = 315
$wsForecast.Cells.Item(60, 2).Value = 316
$wsForecast.Cells.Item(61, 2).Value = 317
$wsForecast.Cells.Item(62, 2).Value = 318
$wsForecast.Cells.Item(63, 2).Value = 319
$wsForecast.Cells.Item(64, 2).Value = 321
$wsForecast.Cells.Item(65, 2).Value = 322
$wsForecast.Cells.Item(66, 2).Value = 323
$wsForecast.Cells.Item(67, 2).Value = 325
$wsForecast.Cells.Item(68, 2).Value = 326
$wsForecast.Cells.Item(69, 2).Value = 327
$wsForecast.Cells.Item(70, 2).Value = 328
$wsForecast.Cells.Item(71, 2).Value = 329
$wsForecast.Cells.Item(72, 2).Value = 330
$wsForecast.Cells.Item(73, 2).Value = 331
$wsForecast.Cells.Item(74, 2).Value = 332
$wsForecast.Cells.Item(75, 2).Value = 333
$wsForecast.Cells.Item(76, 2).Value = 335

# --- Sheet "PO Forecast": re-forecast the tail + extend one more week ---
$wsForecast.Cells.Item(77, 1).Value = 45662.99999999999
$wsForecast.Cells.Item(77, 1).NumberFormat = $wsForecast.Cells.Item(76, 1).NumberFormat
$wsForecast.Cells.Item(77, 2).Value = 344
$wsForecast.Cells.Item(78, 1).Value = 45669.99999999999
$wsForecast.Cells.Item(78, 1).NumberFormat = $wsForecast.Cells.Item(76, 1).NumberFormat
$wsForecast.Cells.Item(78, 2).Value = 345
$wsForecast.Cells.Item(79, 1).Value = 45676.99999999999
$wsForecast.Cells.Item(79, 1).NumberFormat = $wsForecast.Cells.Item(76, 1).NumberFormat
$wsForecast.Cells.Item(79, 2).Value = 346
$wsForecast.Cells.Item(80, 1).Value = 45683.99999999999
$wsForecast.Cells.Item(80, 1).NumberFormat = $wsForecast.Cells.Item(76, 1).NumberFormat
$wsForecast.Cells.Item(80, 2).Value = 347
$wsForecast.Cells.Item(81, 1).Value = 45690.99999999999
$wsForecast.Cells.Item(81, 1).NumberFormat = $wsForecast.Cells.Item(76, 1).NumberFormat
$wsForecast.Cells.Item(81, 2).Value = 348
$wsForecast.Cells.Item(82, 1).Value = 45697.99999999999
$wsForecast.Cells.Item(82, 1).NumberFormat = $wsForecast.Cells.Item(76, 1).NumberFormat
$wsForecast.Cells.Item(82, 2).Value = 349
$wsForecast.Cells.Item(83, 1).Value = 45704.99999999999
$wsForecast.Cells.Item(83, 1).NumberFormat = $wsForecast.Cells.Item(76, 1).NumberFormat
$wsForecast.Cells.Item(83, 2).Value = 350
$wsForecast.Cells.Item(84, 1).Value = 45711.99999999999
$wsForecast.Cells.Item(84, 1).NumberFormat = $wsForecast.Cells.Item(76, 1).NumberFormat
$wsForecast.Cells.Item(84, 2).Value = 351
$wsForecast.Cells.Item(85, 1).Value = 45718.99999999999
$wsForecast.Cells.Item(85, 1).NumberFormat = $wsForecast.Cells.Item(76, 1).NumberFormat
$wsForecast.Cells.Item(85, 2).Value = 352
